$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set for rows 2-5 (header stays in row 1):
#   row 2: Cerveza Corona                       (was row 3)
#   row 3: Ron viejo de caldas (5años) botella   (was row 4)
#   row 4: Aguardiente Amarillo Media            (new product)
#   row 5: Aguardiente Amarillo Botella          (was row 2, re-appended with updated figures)

$rows = @(
    @("Cerveza Corona", "Cervezas", 982, 10000, "31/1/2026", 5000),
    @("Ron viejo de caldas (5años) botella", "Rones", 999, 132000, "3/2/2026", 78000),
    @("Aguardiente Amarillo Media", "Aguardientes", 119, 70000, "6/2/2026", 39000),
    @("Aguardiente Amarillo Botella", "Aguardientes", 1000, 108000, "7/2/2026", 78000)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]

    # Column E holds date-looking text (e.g. "31/1/2026") that must stay a
    # literal string, not get auto-converted into a date serial number.
    # Force the cell to Text format first, assign, then drop back to the
    # default "Normal" style so no lingering number-format is left behind.
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $data[4]
    $eCell.Style = "Normal"

    $ws.Cells.Item($r, 6).Value = $data[5]
}
